$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style_D2 = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.895.78'
$ws.Range('D2').Style = $style_D2
$ws.Range('E2').Value = '  +0.68%  '
$style_D3 = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.633.97'
$ws.Range('D3').Style = $style_D3
$ws.Range('E3').Value = '  +1.70%  '
$ws.Range('E4').Value = '  +0.18%  '
$style_D5 = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.86'
$ws.Range('D5').Style = $style_D5
$ws.Range('E5').Value = '  +0.92%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('E7').Value = '  +0.10%  '
$style_D8 = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '28.75'
$ws.Range('D8').Style = $style_D8
$ws.Range('E8').Value = '  +2.68%  '
$ws.Range('E9').Value = '  +2.29%  '
$style_D10 = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0609'
$ws.Range('D10').Style = $style_D10
$ws.Range('E10').Value = '  +0.70%  '
$style_D11 = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0913'
$ws.Range('D11').Style = $style_D11
$ws.Range('E11').Value = '  +0.36%  '
$style_D12 = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.870.68'
$ws.Range('D12').Style = $style_D12
$ws.Range('E12').Value = '  +1.89%  '
$style_D13 = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.628.71'
$ws.Range('D13').Style = $style_D13
$ws.Range('E13').Value = '  +1.54%  '
$style_D14 = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.566'
$ws.Range('D14').Style = $style_D14
$ws.Range('E14').Value = '  +2.99%  '
$style_D15 = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '9.29'
$ws.Range('D15').Style = $style_D15
$ws.Range('E15').Value = '  +18.26%  '
$style_D16 = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.87'
$ws.Range('D16').Style = $style_D16
$ws.Range('E16').Value = '  +2.87%  '
$style_D17 = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.912.32'
$ws.Range('D17').Style = $style_D17
$ws.Range('E17').Value = '  +0.71%  '
$style_D18 = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '64.31'
$ws.Range('D18').Style = $style_D18
$ws.Range('E18').Value = '  +0.29%  '
$style_D19 = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.46'
$ws.Range('D19').Style = $style_D19
$ws.Range('E19').Value = '  +0.32%  '
$style_D20 = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0700'
$ws.Range('D20').Style = $style_D20
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('E21').Value = '  +0.11%  '
$style_D22 = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.92'
$ws.Range('D22').Style = $style_D22
$ws.Range('E22').Value = '  +5.58%  '
$style_D23 = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.14'
$ws.Range('D23').Style = $style_D23
$ws.Range('E23').Value = '  +2.53%  '
$ws.Range('E24').Value = '  +1.50%  '
$style_D25 = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.00'
$ws.Range('D25').Style = $style_D25
$style_D26 = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.58'
$ws.Range('D26').Style = $style_D26
$ws.Range('E26').Value = '  +0.59%  '
$ws.Range('E27').Value = '  +1.17%  '
$ws.Range('E28').Value = '  +2.40%  '
$ws.Range('E29').Value = '  +0.17%  '
$style_D30 = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0488'
$ws.Range('D30').Style = $style_D30
$ws.Range('E30').Value = '  +1.36%  '
$ws.Range('E31').Value = '  +5.00%  '
$style_D32 = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.37'
$ws.Range('D32').Style = $style_D32
$ws.Range('E32').Value = '  +3.80%  '
$ws.Range('E33').Value = '  -0.17%  '
$style_D34 = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.433.51'
$ws.Range('D34').Style = $style_D34
$ws.Range('E34').Value = '  +0.29%  '
$style_D35 = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.64'
$ws.Range('D35').Style = $style_D35
$ws.Range('E35').Value = '  +4.94%  '
$ws.Range('E36').Value = '  +0.95%  '
$style_D37 = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.83'
$ws.Range('D37').Style = $style_D37
$ws.Range('E37').Value = '  -3.46%  '
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$style_D39 = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0170'
$ws.Range('D39').Style = $style_D39
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$style_D40 = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '76.17'
$ws.Range('D40').Style = $style_D40
$ws.Range('E40').Value = '  +14.45%  '
$style_D41 = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.554'
$ws.Range('D41').Style = $style_D41
$ws.Range('E41').Value = '  +0.77%  '
$style_D42 = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.99'
$ws.Range('D42').Style = $style_D42
$ws.Range('E42').Value = '  +1.76%  '
$style_D43 = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.828'
$ws.Range('D43').Style = $style_D43
$ws.Range('E43').Value = '  +1.29%  '
$style_D44 = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0493'
$ws.Range('D44').Style = $style_D44
$ws.Range('E44').Value = '  -1.32%  '
$style_D45 = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '53.74'
$ws.Range('D45').Style = $style_D45
$ws.Range('E45').Value = '  -6.08%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$style_D46 = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('D46').Style = $style_D46
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$style_D47 = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.01'
$ws.Range('D47').Style = $style_D47
$ws.Range('E47').Value = '  +3.08%  '
$style_D48 = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.777.69'
$ws.Range('D48').Style = $style_D48
$ws.Range('E48').Value = '  +2.14%  '
$style_D49 = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.35'
$ws.Range('D49').Style = $style_D49
$ws.Range('E49').Value = '  +0.09%  '
$style_D50 = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '89.21'
$ws.Range('D50').Style = $style_D50
$ws.Range('E50').Value = '  +2.90%  '
$style_D51 = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0110'
$ws.Range('D51').Style = $style_D51
$ws.Range('E51').Value = '  +4.93%  '
